# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1, styled like the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Values for H2:H34 (0/1 "Save" indicator per row)
$saveValues = @(0,1,0,1,0,0,0,0,0,0,0,0,1,1,0,0,0,1,0,0,0,0,0,0,0,0,1,1,1,1,0,0,1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
